$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2..19) each got the contents of another (old) row for the
# columns D (Fecha), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), N (Unidad de
# comercializacion), O (Origen) and P (Precio $/Kg) -- effectively a
# row-wise shuffle of those columns while A,B,C,E,F,G,H,Q,R stay put.

$cols = @("D","I","J","K","L","M","N","O","P")

# Snapshot the current (pre-edit) values for the columns that move.
$data = @{}
foreach ($r in 2..19) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $data[$r] = $row
}

# Destination row -> source row (source row's old values are written into
# the destination row).
$mapping = @{
    2  = 14
    3  = 19
    4  = 8
    5  = 15
    6  = 6
    7  = 11
    8  = 12
    9  = 13
    10 = 16
    11 = 7
    12 = 17
    13 = 18
    14 = 10
    15 = 9
    16 = 4
    17 = 5
    18 = 2
    19 = 3
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $data[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcData[$c]
    }
}
